$d = $word.ActiveDocument

$replacements = @(
    ,@("2022-12-10 Saturday", "2022-12-11 Sunday")
    ,@("22+76=", "67-38=")
    ,@("79-6=", "29-18=")
    ,@("57-40=", "99-90=")
    ,@("39-5=", "69-34=")
    ,@("29+32=", "81-55=")
    ,@("86-2=", "73-14=")
    ,@("10+14=", "63+34=")
    ,@("58-27=", "37+41=")
    ,@("65-5=", "54-32=")
    ,@("95-14=", "69-37=")
    ,@("66-61=", "20+28=")
    ,@("9+8=", "4+22=")
    ,@("70-10=", "15+33=")
    ,@("77-7=", "72-63=")
    ,@("73-30=", "78-53=")
    ,@("19+5=", "31+66=")
    ,@("24+36=", "5+41=")
    ,@("72-19=", "26+48=")
    ,@("75+10=", "25+9=")
    ,@("89-14=", "23+65=")
    ,@("28+27=", "88-49=")
    ,@("21+43=", "54+13=")
    ,@("85-2=", "43-23=")
    ,@("30+12=", "56-13=")
    ,@("84-54=", "82-44=")
    ,@("55-37=", "5+3=")
    ,@("32+67=", "80-71=")
    ,@("95-70=", "96-78=")
    ,@("99-24=", "54+2=")
    ,@("96-53=", "58-17=")
    ,@("26+69=", "59+6=")
    ,@("89-29=", "29+23=")
    ,@("49+18=", "10+11=")
    ,@("69-10=", "94-7=")
    ,@("71-31=", "3+51=")
    ,@("57-38=", "3+31=")
    ,@("32+32=", "79-61=")
    ,@("34+9=", "57+33=")
    ,@("88-64=", "89-45=")
    ,@("23+33=", "25+37=")
    ,@("81+2=", "60-4=")
    ,@("32-13=", "46-33=")
    ,@("96-56=", "18+34=")
    ,@("64+34=", "95-79=")
    ,@("41-14=", "69-0=")
    ,@("96-2=", "97-7=")
    ,@("91-56=", "50+2=")
    ,@("39-28=", "53+21=")
    ,@("27+68=", "14+33=")
    ,@("28+47=", "75+2=")
    ,@("53-20=", "44+29=")
    ,@("76-47=", "68-15=")
    ,@("75-51=", "94-6=")
    ,@("7+9=", "58-55=")
    ,@("21+56=", "72-0=")
    ,@("30+61=", "8+81=")
    ,@("6+33=", "86-56=")
    ,@("50-13=", "59-7=")
    ,@("97-68=", "70-55=")
    ,@("29+55=", "59-33=")
    ,@("34-16=", "23+31=")
    ,@("73-70=", "81+10=")
    ,@("68+12=", "59+17=")
    ,@("41-30=", "53-40=")
    ,@("47+39=", "72-26=")
    ,@("89-2=", "25+40=")
    ,@("50-19=", "48+27=")
    ,@("0+22=", "42+36=")
    ,@("39+2=", "40-25=")
    ,@("17+22=", "23+14=")
    ,@("17+8=", "22+3=")
    ,@("42-40=", "22+71=")
    ,@("35-11=", "12+40=")
    ,@("77-77=", "4+27=")
    ,@("15+77=", "40+37=")
    ,@("65-61=", "73-64=")
    ,@("69-66=", "18+31=")
    ,@("59+36=", "6+48=")
    ,@("80-56=", "7+18=")
    ,@("1+31=", "63+19=")
    ,@("87-76=", "27-23=")
    ,@("76-58=", "93+6=")
    ,@("86+6=", "88-84=")
    ,@("49+34=", "35+38=")
    ,@("28-15=", "21+46=")
    ,@("96-34=", "1+6=")
    ,@("27+46=", "95-69=")
    ,@("16-15=", "34-22=")
    ,@("27+48=", "10+44=")
    ,@("83-69=", "6+61=")
    ,@("32+4=", "85-12=")
    ,@("49+32=", "70-44=")
    ,@("84-71=", "43+9=")
    ,@("70-42=", "2+88=")
    ,@("9+37=", "56-47=")
    ,@("92-44=", "89-31=")
    ,@("53+36=", "62+3=")
    ,@("55-52=", "17+40=")
    ,@("57-50=", "14+67=")
    ,@("10+51=", "40-26=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Replaced $($replacements.Count) text runs"